$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 6.240107999999999
$ws.Range("H2").Value = 18.720324
$ws.Range("I2").Value = 0.01732230523539376
$ws.Range("J2").Value = 0.01732230523539376
$ws.Range("M2").Value = 0.8793530000000001
$ws.Range("N2").Value = 2.638059
$ws.Range("O2").Value = 0.1744485881486592
$ws.Range("P2").Value = 0.1744485881486592
$ws.Range("Q2").Value = 5.487257690123999
$ws.Range("R2").Value = 49.385319211116
$ws.Range("S2").Value = 0.003021851691794568
$ws.Range("T2").Value = 0.003021851691794569

$ws.Range("G3").Value = 6.240107999999999
$ws.Range("H3").Value = 18.720324
$ws.Range("I3").Value = 0.01732230523539376
$ws.Range("J3").Value = 0.01732230523539376
$ws.Range("O3").Value = 0.6589499911124466
$ws.Range("P3").Value = 0.6589499911124466
$ws.Range("Q3").Value = 20.727186413556
$ws.Range("R3").Value = 186.544677722004
$ws.Range("S3").Value = 0.0114145328809098
$ws.Range("T3").Value = 0.01141453288090981

$ws.Range("G4").Value = 6.240107999999999
$ws.Range("H4").Value = 18.720324
$ws.Range("I4").Value = 0.01732230523539376
$ws.Range("J4").Value = 0.01732230523539376
$ws.Range("O4").Value = 0.1666014207388943
$ws.Range("P4").Value = 0.1666014207388943
$ws.Range("Q4").Value = 5.240426058111999
$ws.Range("R4").Value = 47.16383452300799
$ws.Range("S4").Value = 0.002885920662689387
$ws.Range("T4").Value = 0.002885920662689387

$ws.Range("I5").Value = 0.9592798330716089
$ws.Range("J5").Value = 0.9592798330716091
$ws.Range("M5").Value = 0.8793530000000001
$ws.Range("N5").Value = 2.638059
$ws.Range("O5").Value = 0.1744485881486592
$ws.Range("P5").Value = 0.1744485881486592
$ws.Range("Q5").Value = 303.875007943387
$ws.Range("R5").Value = 2734.875071490483
$ws.Range("S5").Value = 0.1673450125188236
$ws.Range("T5").Value = 0.1673450125188236

$ws.Range("I6").Value = 0.9592798330716089
$ws.Range("J6").Value = 0.9592798330716091
$ws.Range("O6").Value = 0.6589499911124466
$ws.Range("P6").Value = 0.6589499911124466
$ws.Range("S6").Value = 0.6321174374768859
$ws.Range("T6").Value = 0.6321174374768861

$ws.Range("I7").Value = 0.9592798330716089
$ws.Range("J7").Value = 0.9592798330716091
$ws.Range("O7").Value = 0.1666014207388943
$ws.Range("P7").Value = 0.1666014207388943
$ws.Range("S7").Value = 0.1598173830758993
$ws.Range("T7").Value = 0.1598173830758994

$ws.Range("G8").Value = 8.428738666666666
$ws.Range("I8").Value = 0.02339786169299727
$ws.Range("J8").Value = 0.02339786169299728
$ws.Range("M8").Value = 0.8793530000000001
$ws.Range("N8").Value = 2.638059
$ws.Range("O8").Value = 0.1744485881486592
$ws.Range("P8").Value = 0.1744485881486592
$ws.Range("Q8").Value = 7.411836632749333
$ws.Range("R8").Value = 66.70652969474401
$ws.Range("S8").Value = 0.004081723938040969
$ws.Range("T8").Value = 0.004081723938040971

$ws.Range("G9").Value = 8.428738666666666
$ws.Range("I9").Value = 0.02339786169299727
$ws.Range("J9").Value = 0.02339786169299728
$ws.Range("O9").Value = 0.6589499911124466
$ws.Range("P9").Value = 0.6589499911124466
$ws.Range("Q9").Value = 27.99695735637066
$ws.Range("S9").Value = 0.0154180207546508
$ws.Range("T9").Value = 0.01541802075465081

$ws.Range("G10").Value = 8.428738666666666
$ws.Range("I10").Value = 0.02339786169299727
$ws.Range("J10").Value = 0.02339786169299728
$ws.Range("O10").Value = 0.1666014207388943
$ws.Range("P10").Value = 0.1666014207388943
$ws.Range("Q10").Value = 7.078432255630221
$ws.Range("S10").Value = 0.003898117000305495
$ws.Range("T10").Value = 0.003898117000305496
